$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting rows 56..127 down to 57..128
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record's data
$ws.Cells.Item(56, 1).Value = 7
$ws.Cells.Item(56, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(56, 3).Value = "Ñuble"
$ws.Cells.Item(56, 4).Value = 44483
$ws.Cells.Item(56, 5).Value = 16
$ws.Cells.Item(56, 6).Value = 100112017
$ws.Cells.Item(56, 7).Value = "Apio"
$ws.Cells.Item(56, 8).Value = "Americana (o)"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 160
$ws.Cells.Item(56, 11).Value = 8000
$ws.Cells.Item(56, 12).Value = 8500
$ws.Cells.Item(56, 13).Value = 8250
$ws.Cells.Item(56, 14).Value = "$/docena de matas"
$ws.Cells.Item(56, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(56, 16).Value = 1375
$ws.Cells.Item(56, 17).Value = 6
$ws.Cells.Item(56, 18).Value = "Hortaliza"
